$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 new rows at 697, pushing existing rows 697-782 down to 699-784
$ws.Range("A697:A698").EntireRow.Insert()

# New row 697: Pintón, new weekly price observation
$ws.Range("A697").Value = 7
$ws.Range("B697").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C697").Value = "Ñuble"
$ws.Range("D697").Value = 44918
$ws.Range("E697").Value = 16
$ws.Range("F697").Value = "Fruta"
$ws.Range("G697").Value = 100108
$ws.Range("H697").Value = "Tropicales y subtropicales"
$ws.Range("I697").Value = 100108006
$ws.Range("J697").Value = "Plátano"
$ws.Range("K697").Value = "Sin especificar"
$ws.Range("L697").Value = "Pintón"
$ws.Range("M697").Value = 80
$ws.Range("N697").Value = 20000
$ws.Range("O697").Value = 20000
$ws.Range("P697").Value = 20000
$ws.Range("Q697").Value = "$/caja 20 kilos"
$ws.Range("R697").Value = "Ecuador"
$ws.Range("S697").Value = 1000
$ws.Range("T697").Value = 20

# New row 698: Primera Pintón, new weekly price observation
$ws.Range("A698").Value = 7
$ws.Range("B698").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C698").Value = "Ñuble"
$ws.Range("D698").Value = 44918
$ws.Range("E698").Value = 16
$ws.Range("F698").Value = "Fruta"
$ws.Range("G698").Value = 100108
$ws.Range("H698").Value = "Tropicales y subtropicales"
$ws.Range("I698").Value = 100108006
$ws.Range("J698").Value = "Plátano"
$ws.Range("K698").Value = "Sin especificar"
$ws.Range("L698").Value = "Primera Pintón"
$ws.Range("M698").Value = 160
$ws.Range("N698").Value = 21000
$ws.Range("O698").Value = 25000
$ws.Range("P698").Value = 23000
$ws.Range("Q698").Value = "$/caja 20 kilos"
$ws.Range("R698").Value = "Ecuador"
$ws.Range("S698").Value = 1150
$ws.Range("T698").Value = 20
